# Update countries & provincias Spain
# Applies the daily COVID-data refresh to the "Pais" sheet:
#  - updates the "Datos actualizados..." timestamp
#  - updates numeric stats for several countries (some rows simply got
#    new totals, others were re-sorted causing two adjacent rows to swap
#    which country/data they hold)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 23:08"

# --- Estados Unidos (row 4): totals refreshed, no reorder ---
$ws.Range("B4").Value = 1256039
$ws.Range("C4").Value = 18406
$ws.Range("E4").Value = 976731
$ws.Range("F4").Value = 15851
$ws.Range("G4").Value = 1822
$ws.Range("H4").Value = 74093

# --- Alemania (row 9): totals refreshed, no reorder ---
$ws.Range("B9").Value = 167817
$ws.Range("C9").Value = 810
$ws.Range("D9").Value = 137696
$ws.Range("E9").Value = 22896
$ws.Range("G9").Value = 232
$ws.Range("H9").Value = 7225

# --- Canada (row 15): totals refreshed, no reorder ---
$ws.Range("B15").Value = 63403
$ws.Range("C15").Value = 1357
$ws.Range("D15").Value = 27816
$ws.Range("E15").Value = 31364

# --- Rows 84/85: Nueva Zelanda <-> Costa de Marfil swap places ---
$ws.Range("A84").Value = "Costa de Marfil"
$ws.Range("B84").Value = 1516
$ws.Range("C84").Value = 52
$ws.Range("D84").Value = 721
$ws.Range("E84").Value = 777
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 18

$ws.Range("A85").Value = "Nueva Zelanda"
$ws.Range("B85").Value = 1488
$ws.Range("C85").Value = 2
$ws.Range("D85").Value = 1316
$ws.Range("E85").Value = 151
$ws.Range("F85").Value = 2
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 21

# --- Rows 102/103: Guatemala <-> Niger swap places ---
$ws.Range("A102").Value = "Niger"
$ws.Range("B102").Value = 770
$ws.Range("C102").Value = 7
$ws.Range("D102").Value = 561
$ws.Range("E102").Value = 171
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 38

$ws.Range("A103").Value = "Guatemala"
$ws.Range("B103").Value = 763
$ws.Range("C103").Value = 33
$ws.Range("D103").Value = 79
$ws.Range("E103").Value = 665
$ws.Range("F103").Value = 5
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 19

# --- Maldivas (row 113): totals refreshed, no reorder ---
$ws.Range("E113").Value = 595
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 2

# --- Rows 158/159: Nepal <-> Uganda swap places ---
$ws.Range("A158").Value = "Uganda"
$ws.Range("B158").Value = 100
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 55
$ws.Range("E158").Value = 45

$ws.Range("A159").Value = "Nepal"
$ws.Range("B159").Value = 99
$ws.Range("C159").Value = 17
$ws.Range("D159").Value = 22
$ws.Range("E159").Value = 77

# --- Rows 205/206: Seychelles <-> Montserrat swap places ---
$ws.Range("A205").Value = "Montserrat"
$ws.Range("D205").Value = 7
$ws.Range("F205").Value = 1
$ws.Range("H205").Value = 1

$ws.Range("A206").Value = "Seychelles"
$ws.Range("D206").Value = 8
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0
